$wb = $excel.ActiveWorkbook

# --- ALC sheet: row 138 value updates ---
$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H138").Value = 14401.4
$wsALC.Range("J138").Value = 594
$wsALC.Range("L138").Value = 1782
$wsALC.Range("N138").Value = -12062

# --- CRP sheet: add H-N values for rows previously blank ---
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Range("H129").Value = 0
$wsCRP.Range("I129").Value = 0
$wsCRP.Range("J129").Value = 0
$wsCRP.Range("K129").Value = 0
$wsCRP.Range("L129").Value = 0
$wsCRP.Range("H130").Value = 0
$wsCRP.Range("I130").Value = 0
$wsCRP.Range("J130").Value = 0
$wsCRP.Range("K130").Value = 0
$wsCRP.Range("L130").Value = 0
$wsCRP.Range("H131").Value = 0
$wsCRP.Range("I131").Value = 0
$wsCRP.Range("J131").Value = 0
$wsCRP.Range("K131").Value = 0
$wsCRP.Range("L131").Value = 0
$wsCRP.Range("H132").Value = 5246.75
$wsCRP.Range("I132").Value = 4996
$wsCRP.Range("J132").Value = 5999
$wsCRP.Range("K132").Value = 14988
$wsCRP.Range("L132").Value = 17997
$wsCRP.Range("M132").Value = -12458
$wsCRP.Range("N132").Value = -23057
$wsCRP.Range("H133").Value = 0
$wsCRP.Range("I133").Value = 0
$wsCRP.Range("J133").Value = 0
$wsCRP.Range("K133").Value = 0
$wsCRP.Range("L133").Value = 0
$wsCRP.Range("H134").Value = 3000
$wsCRP.Range("I134").Value = 2375
$wsCRP.Range("J134").Value = 4250
$wsCRP.Range("K134").Value = 7125
$wsCRP.Range("L134").Value = 12750
$wsCRP.Range("M134").Value = -4590
$wsCRP.Range("N134").Value = -17820
$wsCRP.Range("H135").Value = 70000
$wsCRP.Range("I135").Value = 0
$wsCRP.Range("J135").Value = 70000
$wsCRP.Range("K135").Value = 0
$wsCRP.Range("L135").Value = 70000
$wsCRP.Range("N135").Value = -80140
$wsCRP.Range("H137").Value = 100000
$wsCRP.Range("I137").Value = 0
$wsCRP.Range("J137").Value = 100000
$wsCRP.Range("K137").Value = 0
$wsCRP.Range("L137").Value = 100000
$wsCRP.Range("N137").Value = -110200
$wsCRP.Range("H138").Value = 30000
$wsCRP.Range("I138").Value = 0
$wsCRP.Range("J138").Value = 30000
$wsCRP.Range("K138").Value = 0
$wsCRP.Range("L138").Value = 30000
$wsCRP.Range("N138").Value = -40280
$wsCRP.Range("H139").Value = 0
$wsCRP.Range("I139").Value = 0
$wsCRP.Range("J139").Value = 0
$wsCRP.Range("K139").Value = 0
$wsCRP.Range("L139").Value = 0
$wsCRP.Range("H140").Value = 0
$wsCRP.Range("I140").Value = 0
$wsCRP.Range("J140").Value = 0
$wsCRP.Range("K140").Value = 0
$wsCRP.Range("L140").Value = 0
$wsCRP.Range("H141").Value = 0
$wsCRP.Range("I141").Value = 0
$wsCRP.Range("J141").Value = 0
$wsCRP.Range("K141").Value = 0
$wsCRP.Range("L141").Value = 0

# --- CUL sheet: clear H-N for rows (135 untouched) ---
$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("H120:N120").ClearContents()
$wsCUL.Range("H121:N121").ClearContents()
$wsCUL.Range("H122:N122").ClearContents()
$wsCUL.Range("H123:N123").ClearContents()
$wsCUL.Range("H124:N124").ClearContents()
$wsCUL.Range("H125:N125").ClearContents()
$wsCUL.Range("H126:N126").ClearContents()
$wsCUL.Range("H127:N127").ClearContents()
$wsCUL.Range("H128:N128").ClearContents()
$wsCUL.Range("H129:N129").ClearContents()
$wsCUL.Range("H130:N130").ClearContents()
$wsCUL.Range("H131:N131").ClearContents()
$wsCUL.Range("H132:N132").ClearContents()
$wsCUL.Range("H133:N133").ClearContents()
$wsCUL.Range("H134:N134").ClearContents()
$wsCUL.Range("H136:N136").ClearContents()
$wsCUL.Range("H137:N137").ClearContents()
$wsCUL.Range("H138:N138").ClearContents()
$wsCUL.Range("H139:N139").ClearContents()
$wsCUL.Range("H140:N140").ClearContents()
$wsCUL.Range("H141:N141").ClearContents()

# --- GSM sheet: update rows 80 & 83, add H-N for rows 125-141 ---
$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("H125").Value = 0
$wsGSM.Range("I125").Value = 0
$wsGSM.Range("J125").Value = 0
$wsGSM.Range("K125").Value = 0
$wsGSM.Range("L125").Value = 0
$wsGSM.Range("H126").Value = 0
$wsGSM.Range("I126").Value = 0
$wsGSM.Range("J126").Value = 0
$wsGSM.Range("K126").Value = 0
$wsGSM.Range("L126").Value = 0
$wsGSM.Range("H127").Value = 0
$wsGSM.Range("I127").Value = 0
$wsGSM.Range("J127").Value = 0
$wsGSM.Range("K127").Value = 0
$wsGSM.Range("L127").Value = 0
$wsGSM.Range("H128").Value = 80000
$wsGSM.Range("I128").Value = 0
$wsGSM.Range("J128").Value = 80000
$wsGSM.Range("K128").Value = 0
$wsGSM.Range("L128").Value = 80000
$wsGSM.Range("N128").Value = -89960
$wsGSM.Range("H129").Value = 0
$wsGSM.Range("I129").Value = 0
$wsGSM.Range("J129").Value = 0
$wsGSM.Range("K129").Value = 0
$wsGSM.Range("L129").Value = 0
$wsGSM.Range("H130").Value = 0
$wsGSM.Range("I130").Value = 0
$wsGSM.Range("J130").Value = 0
$wsGSM.Range("K130").Value = 0
$wsGSM.Range("L130").Value = 0
$wsGSM.Range("H131").Value = 0
$wsGSM.Range("I131").Value = 0
$wsGSM.Range("J131").Value = 0
$wsGSM.Range("K131").Value = 0
$wsGSM.Range("L131").Value = 0
$wsGSM.Range("H132").Value = 4054.5
$wsGSM.Range("I132").Value = 3960.05
$wsGSM.Range("J132").Value = 4999
$wsGSM.Range("K132").Value = 11880.15
$wsGSM.Range("L132").Value = 14997
$wsGSM.Range("M132").Value = -9350.150000000001
$wsGSM.Range("N132").Value = -20057
$wsGSM.Range("H133").Value = 0
$wsGSM.Range("I133").Value = 0
$wsGSM.Range("J133").Value = 0
$wsGSM.Range("K133").Value = 0
$wsGSM.Range("L133").Value = 0
$wsGSM.Range("H134").Value = 0
$wsGSM.Range("I134").Value = 0
$wsGSM.Range("J134").Value = 0
$wsGSM.Range("K134").Value = 0
$wsGSM.Range("L134").Value = 0
$wsGSM.Range("H135").Value = 0
$wsGSM.Range("I135").Value = 0
$wsGSM.Range("J135").Value = 0
$wsGSM.Range("K135").Value = 0
$wsGSM.Range("L135").Value = 0
$wsGSM.Range("H136").Value = 0
$wsGSM.Range("I136").Value = 0
$wsGSM.Range("J136").Value = 0
$wsGSM.Range("K136").Value = 0
$wsGSM.Range("L136").Value = 0
$wsGSM.Range("H137").Value = 0
$wsGSM.Range("I137").Value = 0
$wsGSM.Range("J137").Value = 0
$wsGSM.Range("K137").Value = 0
$wsGSM.Range("L137").Value = 0
$wsGSM.Range("H138").Value = 150000
$wsGSM.Range("I138").Value = 0
$wsGSM.Range("J138").Value = 150000
$wsGSM.Range("K138").Value = 0
$wsGSM.Range("L138").Value = 150000
$wsGSM.Range("N138").Value = -160280
$wsGSM.Range("H139").Value = 0
$wsGSM.Range("I139").Value = 0
$wsGSM.Range("J139").Value = 0
$wsGSM.Range("K139").Value = 0
$wsGSM.Range("L139").Value = 0
$wsGSM.Range("H140").Value = 0
$wsGSM.Range("I140").Value = 0
$wsGSM.Range("J140").Value = 0
$wsGSM.Range("K140").Value = 0
$wsGSM.Range("L140").Value = 0
$wsGSM.Range("H141").Value = 0
$wsGSM.Range("I141").Value = 0
$wsGSM.Range("J141").Value = 0
$wsGSM.Range("K141").Value = 0
$wsGSM.Range("L141").Value = 0

# GSM rows 80 and 83: update values and clear N column (previously had N, now removed)
$wsGSM.Range("H80").Value = 3000
$wsGSM.Range("I80").Value = 3000
$wsGSM.Range("J80").Value = 0
$wsGSM.Range("K80").Value = 3000
$wsGSM.Range("L80").Value = 0
$wsGSM.Range("M80").Value = -2002
$wsGSM.Range("N80").ClearContents()
$wsGSM.Range("H83").Value = 3000
$wsGSM.Range("I83").Value = 3000
$wsGSM.Range("J83").Value = 0
$wsGSM.Range("K83").Value = 15000
$wsGSM.Range("L83").Value = 0
$wsGSM.Range("M83").Value = -10008
$wsGSM.Range("N83").ClearContents()
